$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dbo_detection")

# Insert a new row at row 2 (pushing existing rows 2-5 down to 3-6)
$ws.Rows.Item(2).Insert()

# New row 2: detection_id=0, detection_name="Not available", previous_name="NA"
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 3).Value = "NA"
$ws.Cells.Item(2, 2).Value = "Not available"

# Renumber detection_id column for subsequent rows (now rows 3-6)
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4

# Update the defined name range to extend through the new row
$wb.Names.Item("dbo_biogroup").RefersTo = "=dbo_detection!`$A`$1:`$B`$6"

# Update the active selection to match the new layout
$ws.Range("B5").Select()
